# "Generate Report for Handoff"
#
# The bc4e98ba-a8c9-4cec-b50d-47e6807a0348.md file just finished being
# handed off for localization, so its status flips from "In Translation"
# to "Ready for handoff" (with a fresh handoff timestamp and "mt" priority)
# on every sheet that tracks it: the per-language "zh-cn"/"de-de" sheets,
# and the roll-up "Overview" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 is the bc4e98ba...md file
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = "Ready for handoff"     # Status
$wsZh.Range("E3").Value = "mt"                     # Priority
$wsZh.Range("H3").Value = "2016-09-02 08:16:23"    # Latest Handoff Datetime
$wsZh.Columns.Item(3).ColumnWidth = 16.33

# ---------------------------------------------------------------------
# de-de sheet: row 3 is the bc4e98ba...md file
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = "Ready for handoff"     # Status
$wsDe.Range("E3").Value = "mt"                     # Priority
$wsDe.Range("H3").Value = "2016-09-02 08:16:27"    # Latest Handoff Datetime
$wsDe.Columns.Item(3).ColumnWidth = 16.33

# ---------------------------------------------------------------------
# Overview sheet: row 3 is the bc4e98ba...md file; columns E (zh-cn) and
# F (de-de) hold the per-language status, column G the latest generate date
# ---------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")
$wsOv.Range("E3").Value = "Ready for handoff"
$wsOv.Range("F3").Value = "Ready for handoff"
$wsOv.Range("G3").Value = "2016-09-02 08:16:27"
$wsOv.Columns.Item(5).ColumnWidth = 16.33
$wsOv.Columns.Item(6).ColumnWidth = 16.33
